$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1994047619047619
$ws.Range("C2").Value = 0.5446428571428571
$ws.Range("J2").Value = 0.0119047619047619
$ws.Range("P2").Value = 0.130952380952381
$ws.Range("S2").Value = 0.1130952380952381
$ws.Range("B3").Value = 0.01570680628272251
$ws.Range("C3").Value = 0.03664921465968586
$ws.Range("J3").Value = 0.01570680628272251
$ws.Range("P3").Value = 0.7539267015706806
$ws.Range("S3").Value = 0.1780104712041885
$ws.Range("J4").Value = 0.0576923076923077
$ws.Range("P4").Value = 0.6153846153846154
$ws.Range("S4").Value = 0.3269230769230769
$ws.Range("B6").Value = 0.04784688995215311
$ws.Range("E6").Value = 0.004784688995215311
$ws.Range("F6").Value = 0.03827751196172249
$ws.Range("J6").Value = 0.2966507177033493
$ws.Range("O6").Value = 0.01913875598086124
$ws.Range("Q6").Value = 0.1674641148325359
$ws.Range("R6").Value = 0.07177033492822966
$ws.Range("S6").Value = 0.354066985645933
$ws.Range("B7").Value = 0.1256038647342995
$ws.Range("D7").Value = 0.01932367149758454
$ws.Range("E7").Value = 0.004830917874396135
$ws.Range("F7").Value = 0.04830917874396135
$ws.Range("J7").Value = 0.1594202898550725
$ws.Range("O7").Value = 0.02415458937198068
$ws.Range("Q7").Value = 0.1400966183574879
$ws.Range("R7").Value = 0.07246376811594203
$ws.Range("S7").Value = 0.4057971014492754
$ws.Range("B8").Value = 0.1046228710462287
$ws.Range("D8").Value = 0.0170316301703163
$ws.Range("E8").Value = 0.0024330900243309
$ws.Range("F8").Value = 0.04379562043795621
$ws.Range("J8").Value = 0.09002433090024331
$ws.Range("O8").Value = 0.0340632603406326
$ws.Range("Q8").Value = 0.1557177615571776
$ws.Range("R8").Value = 0.09002433090024331
$ws.Range("S8").Value = 0.462287104622871
$ws.Range("B9").Value = 0.1052631578947368
$ws.Range("D9").Value = 0.005263157894736842
$ws.Range("E9").Value = 0.01052631578947368
$ws.Range("F9").Value = 0.05789473684210526
$ws.Range("J9").Value = 0.1210526315789474
$ws.Range("O9").Value = 0.02105263157894737
$ws.Range("Q9").Value = 0.1210526315789474
$ws.Range("R9").Value = 0.07368421052631578
$ws.Range("S9").Value = 0.4842105263157895
$ws.Range("B10").Value = 0.1345076060848679
$ws.Range("D10").Value = 0.03282626100880705
$ws.Range("E10").Value = 0.001601281024819856
$ws.Range("F10").Value = 0.07205764611689351
$ws.Range("J10").Value = 0.0976781425140112
$ws.Range("O10").Value = 0.02562049639711769
$ws.Range("Q10").Value = 0.1929543634907926
$ws.Range("R10").Value = 0.06805444355484387
$ws.Range("S10").Value = 0.3746997598078463
$ws.Range("G11").Value = 0.1725239616613418
$ws.Range("J11").Value = 0.08626198083067092
$ws.Range("K11").Value = 0.194888178913738
$ws.Range("L11").Value = 0.5335463258785943
$ws.Range("S11").Value = 0.01277955271565495
$ws.Range("G12").Value = 0.7251461988304093
$ws.Range("J12").Value = 0.1871345029239766
$ws.Range("K12").Value = 0.005847953216374269
$ws.Range("L12").Value = 0.02339181286549707
$ws.Range("S12").Value = 0.05847953216374269
$ws.Range("G13").Value = 0.6122448979591837
$ws.Range("J13").Value = 0.3469387755102041
$ws.Range("S13").Value = 0.04081632653061224
$ws.Range("F15").Value = 0.01976284584980237
$ws.Range("H15").Value = 0.1343873517786561
$ws.Range("I15").Value = 0.05928853754940711
$ws.Range("J15").Value = 0.3517786561264822
$ws.Range("K15").Value = 0.06324110671936758
$ws.Range("M15").Value = 0.007905138339920948
$ws.Range("O15").Value = 0.06719367588932806
$ws.Range("S15").Value = 0.2964426877470356
$ws.Range("F16").Value = 0.01869158878504673
$ws.Range("H16").Value = 0.1214953271028037
$ws.Range("I16").Value = 0.102803738317757
$ws.Range("J16").Value = 0.411214953271028
$ws.Range("K16").Value = 0.102803738317757
$ws.Range("M16").Value = 0.04205607476635514
$ws.Range("O16").Value = 0.06542056074766354
$ws.Range("S16").Value = 0.1355140186915888
$ws.Range("F17").Value = 0.02051282051282051
$ws.Range("H17").Value = 0.1743589743589744
$ws.Range("I17").Value = 0.1025641025641026
$ws.Range("J17").Value = 0.3743589743589744
$ws.Range("K17").Value = 0.09743589743589744
$ws.Range("M17").Value = 0.02564102564102564
$ws.Range("N17").Value = 0.002564102564102564
$ws.Range("O17").Value = 0.06923076923076923
$ws.Range("S17").Value = 0.1333333333333333
$ws.Range("H18").Value = 0.1325301204819277
$ws.Range("I18").Value = 0.09036144578313253
$ws.Range("J18").Value = 0.4819277108433735
$ws.Range("K18").Value = 0.1325301204819277
$ws.Range("M18").Value = 0.006024096385542169
$ws.Range("O18").Value = 0.05421686746987952
$ws.Range("S18").Value = 0.1024096385542169
$ws.Range("F19").Value = 0.021986353297953
$ws.Range("H19").Value = 0.2031842304776346
$ws.Range("I19").Value = 0.07505686125852919
$ws.Range("J19").Value = 0.3805913570887036
$ws.Range("K19").Value = 0.1137225170583776
$ws.Range("M19").Value = 0.02122820318423048
$ws.Range("N19").Value = 0.0007581501137225171
$ws.Range("O19").Value = 0.07505686125852919
$ws.Range("S19").Value = 0.1084154662623199
